# Code Version 29052020 1333
#
# Adds two new lookup rows (PCI-DSS embossing / non-embossing UNC share
# folders) to the bottom of the PostBatchScriptFileLocation lookup sheet,
# with matching hyperlinks, and widens column B (Activity_Name) to fit the
# new, longer activity names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append rows 41 & 42, cloning formatting from row 40 (last existing row) ---
$ws.Range("A40:F40").Copy()
$ws.Range("A41:F41").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A42:F42").PasteSpecial(-4122)   # xlPasteFormats

# Fill in Path (column C) before Activity_Name (column B) for both rows,
# row 41 then row 42 - matches the shared-string insertion order of the
# original edit.
$ws.Cells.Item(41, 3).Value = "\\172.16.17.183\pmu\emboss_files\TSYSCREDIT\embossing"
$ws.Cells.Item(42, 3).Value = "\\172.16.17.183\pmu\emboss_files\TSYSCREDIT\non_embossing"
$ws.Cells.Item(42, 2).Value = "fd_NonEmbossingPcidssFolder"
$ws.Cells.Item(41, 2).Value = "fd_EmbossingPcidssFolder"

# Remaining columns: S.No, IsFolder, Server
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 4).Value = "Yes"
$ws.Cells.Item(41, 6).Value = "Network"

$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 4).Value = "Yes"
$ws.Cells.Item(42, 6).Value = "Network"

# --- Hyperlink the two new Path cells to their UNC share ---
$ws.Hyperlinks.Add($ws.Cells.Item(41, 3), "file:///\\172.16.17.183\pmu\emboss_files\TSYSCREDIT\embossing")
$ws.Hyperlinks.Add($ws.Cells.Item(42, 3), "file:///\\172.16.17.183\pmu\emboss_files\TSYSCREDIT\non_embossing")

# --- Widen column B (Activity_Name) so the new, longer names fit ---
$ws.Columns.Item(2).ColumnWidth = 25.7109375
